# "Fin de Sprint 1" — journal de travail update
#
# Summary of the edit (derived from the canonical OOXML diff):
#  - H11 text corrected: "Dossier de projet: Analyse / Conception: Concept"
#       -> "Dossier de projet: Concept"
#  - F28 / F36 category corrected: "Analyse" -> "Conception"
#  - Row 11 loses its manual 30pt row height (content now fits on one line)
#  - Row 40 finished (D40/E40 filled in) and four new journal rows (41-44)
#    added for the rest of Sprint 1 (Risques techniques, Analyse
#    concurentielle, Base de données / Diagramme MCD)
#  - The now-unused placeholder rows 45-46 are cleared back to the same
#    "empty" shape as the rows below them
#  - Two trailing blank rows (150, 151) at the very end of the sheet are
#    removed
#  - The worksheet's saved scroll position / selection is reset (no more
#    topLeftCell="A13" / selection on H41)
#  - "Rapport" sheet gets an explicit paper size / orientation page setup

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- Row 11: fix the title text, then let Excel re-fit the row height ----
$ws.Range("H11").Value2 = "Dossier de projet: Concept"
$ws.Rows.Item(11).AutoFit()

# --- Rows 28 & 36: these entries were mis-categorised as "Analyse" -------
$ws.Range("F28").Value2 = "Conception"
$ws.Range("F36").Value2 = "Conception"

# --- Row 40: close out the "Stratégie de test" entry ---------------------
$ws.Range("D40").Value2 = 0.39930555555555558
$ws.Range("E40").Formula = "=D40-C40"

# --- Rows 41-44: remainder of Sprint 1 ------------------------------------
$ws.Range("A41").Value2 = 44323
$ws.Range("B41").Value2 = 1
$ws.Range("C41").Value2 = 0.40972222222222227
$ws.Range("D41").Value2 = 0.42708333333333331
$ws.Range("E41").Formula = "=D41-C41"
$ws.Range("F41").Value2 = "Analyse"
$ws.Range("G41").Value2 = "Documentation"
$ws.Range("H41").Value2 = "Dossier de projet: Stratégie de test"

$ws.Range("A42").Value2 = 44323
$ws.Range("B42").Value2 = 1
$ws.Range("C42").Value2 = 0.42708333333333331
$ws.Range("D42").Value2 = 0.44097222222222227
$ws.Range("E42").Formula = "=D42-C42"
$ws.Range("F42").Value2 = "Analyse"
$ws.Range("G42").Value2 = "Documentation"
$ws.Range("H42").Value2 = "Dossier de projet: Risques techniques"

$ws.Range("A43").Value2 = 44323
$ws.Range("B43").Value2 = 1
$ws.Range("C43").Value2 = 0.44097222222222227
$ws.Range("D43").Value2 = 0.46249999999999997
$ws.Range("E43").Formula = "=D43-C43"
$ws.Range("F43").Value2 = "Analyse"
$ws.Range("G43").Value2 = "Documentation"
$ws.Range("H43").Value2 = "Dossier de projet: Analyse concurentielle"

$ws.Range("A44").Value2 = 44323
$ws.Range("B44").Value2 = 1
$ws.Range("C44").Value2 = 0.46249999999999997
$ws.Range("D44").Value2 = 0.4826388888888889
$ws.Range("E44").Formula = "=D44-C44"
$ws.Range("F44").Value2 = "Conception"
$ws.Range("G44").Value2 = "Base de données"
$ws.Range("H44").Value2 = "Diagramme MCD"

# --- Rows 45-46: these were still blank template rows, clear them back to
#     the plain "unused" shape used by the rows that follow them ----------
$ws.Range("A45:D46").Clear()
$ws.Range("G45:L46").Clear()

# --- Drop the two trailing blank rows at the bottom of the sheet ---------
$ws.Rows.Item(150).Delete()
$ws.Rows.Item(150).Delete()

# --- Reset the saved view: no more stale scroll position / selection -----
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# --- "Rapport" sheet: give it an explicit page setup ----------------------
$ws2 = $wb.Worksheets.Item("Rapport")
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1
